$wb = $excel.ActiveWorkbook

# Rename the six "* World Names" / "Void Born Names" sheets, dropping the
# trailing " Names" to match the renamed tabs.
$wb.Worksheets.Item("Death World Names").Name        = "Death World"
$wb.Worksheets.Item("Noble Born World Names").Name    = "Noble Born World"
$wb.Worksheets.Item("Imperial World Names").Name      = "Imperial World"
$wb.Worksheets.Item("Hive World Names").Name          = "Hive World"
$wb.Worksheets.Item("Forge World Names").Name         = "Forge World"
$wb.Worksheets.Item("Void Born Names").Name           = "Void Born"

$wsDeathWorld = $wb.Worksheets.Item("Death World")
$wsVoidBorn   = $wb.Worksheets.Item("Void Born")

# "Death World" loses tab selection and its prior B1:F133 selection,
# landing on A15 instead.
$wsDeathWorld.Activate()
$wsDeathWorld.Range("A15").Select()

# "Void Born" becomes the selected/active tab, with the cursor on M32.
$wsVoidBorn.Activate()
$wsVoidBorn.Range("M32").Select()
